$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text formatting for numeric-looking price values so they retain exact literal formatting
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Apply updated cell values
$ws.Range("D2").Value = "67.867.32"
$ws.Range("E2").Value = "  +0.07%  "
$ws.Range("D3").Value = "3.737.15"
$ws.Range("E3").Value = "  -2.03%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "594.12"
$ws.Range("E5").Value = "  -0.85%  "
$ws.Range("D6").Value = "166.37"
$ws.Range("E6").Value = "  -1.27%  "
$ws.Range("D7").Value = "3.736.70"
$ws.Range("E7").Value = "  -1.95%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("D9").Value = "0.520"
$ws.Range("E9").Value = "  -1.67%  "
$ws.Range("D10").Value = "0.159"
$ws.Range("E10").Value = "  -3.23%  "
$ws.Range("D11").Value = "6.47"
$ws.Range("E11").Value = "  -0.31%  "
$ws.Range("E12").Value = "  -2.59%  "
$ws.Range("E13").Value = "  -4.37%  "
$ws.Range("D14").Value = "36.66"
$ws.Range("E14").Value = "  -0.82%  "
$ws.Range("D15").Value = "4.368.66"
$ws.Range("E15").Value = "  -2.14%  "
$ws.Range("D16").Value = "3.739.00"
$ws.Range("E16").Value = "  -2.24%  "
$ws.Range("D17").Value = "67.890.68"
$ws.Range("E17").Value = "  +0.01%  "
$ws.Range("D18").Value = "18.12"
$ws.Range("E18").Value = "  -2.99%  "
$ws.Range("E19").Value = "  -4.45%  "
$ws.Range("D20").Value = "0.111"
$ws.Range("E20").Value = "  -0.45%  "
$ws.Range("D21").Value = "10.76"
$ws.Range("E21").Value = "  -0.39%  "
$ws.Range("D22").Value = "466.90"
$ws.Range("E22").Value = "  -0.30%  "
$ws.Range("D23").Value = "0.701"
$ws.Range("E23").Value = "  -4.83%  "
$ws.Range("D24").Value = "83.13"
$ws.Range("E24").Value = "  -0.53%  "
$ws.Range("D25").Value = "0.0000137"
$ws.Range("E25").Value = "  -9.20%  "
$ws.Range("D26").Value = "2.22"
$ws.Range("E26").Value = "  -2.70%  "
$ws.Range("D27").Value = "12.08"
$ws.Range("E27").Value = "  -0.76%  "
$ws.Range("D28").Value = "10.14"
$ws.Range("E28").Value = "  -1.83%  "
$ws.Range("E29").Value = "  +0.02%  "
$ws.Range("D30").Value = "3.888.37"
$ws.Range("E30").Value = "  -1.97%  "
$ws.Range("D31").Value = "2.78"
$ws.Range("E31").Value = "  -4.71%  "
$ws.Range("D32").Value = "7.38"
$ws.Range("E32").Value = "  -4.48%  "
$ws.Range("D33").Value = "2.24"
$ws.Range("E33").Value = "  -1.63%  "
$ws.Range("D34").Value = "29.78"
$ws.Range("E34").Value = "  -2.79%  "
$ws.Range("D35").Value = "9.12"
$ws.Range("E35").Value = "  -1.72%  "
$ws.Range("D36").Value = "0.998"
$ws.Range("D37").Value = "3.694.95"
$ws.Range("E37").Value = "  -2.37%  "
$ws.Range("E38").Value = "  -3.90%  "
$ws.Range("D39").Value = "3.42"
$ws.Range("E39").Value = "  -10.73%  "
$ws.Range("E40").Value = "  -1.36%  "
$ws.Range("D41").Value = "0.990"
$ws.Range("E41").Value = "  -2.38%  "
$ws.Range("D42").Value = "5.78"
$ws.Range("E42").Value = "  -3.00%  "
$ws.Range("D43").Value = "1.00"
$ws.Range("E43").Value = "  -0.02%  "
$ws.Range("D45").Value = "0.307"
$ws.Range("E45").Value = "  -3.34%  "
$ws.Range("D46").Value = "8.58"
$ws.Range("E46").Value = "  -1.98%  "
$ws.Range("E47").Value = "  -2.64%  "
$ws.Range("E48").Value = "  -2.27%  "
$ws.Range("D49").Value = "394.62"
$ws.Range("E49").Value = "  -3.01%  "
$ws.Range("D50").Value = "145.07"
$ws.Range("E50").Value = "  +1.46%  "
$ws.Range("D51").Value = "25.62"
$ws.Range("E51").Value = "  +0.65%  "
